$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.688.65"
$ws.Range("E2").Value = "  -2.18%  "

# Row 3
$ws.Range("D3").Value = "1.804.85"
$ws.Range("E3").Value = "  -2.51%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.91"
$ws.Range("E5").Value = "  +0.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.603"
$ws.Range("E6").Value = "  -1.00%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "38.99"
$ws.Range("E8").Value = "  -6.67%  "

# Row 9
$ws.Range("E9").Value = "  +3.25%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0678"
$ws.Range("E10").Value = "  -2.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0991"
$ws.Range("E11").Value = "  -2.04%  "

# Row 12
$ws.Range("D12").Value = "2.065.19"
$ws.Range("E12").Value = "  -2.57%  "

# Row 13
$ws.Range("D13").Value = "1.794.25"
$ws.Range("E13").Value = "  -2.85%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.658"
$ws.Range("E14").Value = "  -1.88%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.85"
$ws.Range("E15").Value = "  -5.13%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.54"
$ws.Range("E16").Value = "  -3.16%  "

# Row 17
$ws.Range("D17").Value = "34.666.59"
$ws.Range("E17").Value = "  -2.22%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.35"
$ws.Range("E18").Value = "  -0.78%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0780"
$ws.Range("E19").Value = "  -2.42%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.97"
$ws.Range("E20").Value = "  -3.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.70"
$ws.Range("E21").Value = "  -3.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.63"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("E23").Value = "  +0.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").Value = "  +1.68%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.55"
$ws.Range("E25").Value = "  +2.01%  "

# Row 26
$ws.Range("E26").Value = "  -2.77%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.10"
$ws.Range("E27").Value = "  -3.61%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.119"
$ws.Range("E28").Value = "  -2.61%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.51"
$ws.Range("E29").Value = "  +8.74%  "

# Row 30
$ws.Range("E30").Value = "  +0.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.98"
$ws.Range("E31").Value = "  +1.44%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0542"
$ws.Range("E32").Value = "  -0.22%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.93"
$ws.Range("E33").Value = "  -3.22%  "

# Row 34
$ws.Range("E34").Value = "  +15.03%  "

# Row 35
$ws.Range("E35").Value = "  -5.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.688"
$ws.Range("E36").Value = "  +0.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "90.70"
$ws.Range("E37").Value = "  -8.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.32"
$ws.Range("E38").Value = "  +4.61%  "

# Row 39
$ws.Range("D39").Value = "1.307.79"
$ws.Range("E39").Value = "  -3.75%  "

# Row 40
$ws.Range("E40").Value = "  -2.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.45"
$ws.Range("E41").Value = "  -0.89%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.954"
$ws.Range("E42").Value = "  -4.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.14"
$ws.Range("E43").Value = "  -3.28%  "

# Row 44
$ws.Range("E44").Value = "  -10.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.65"
$ws.Range("E45").Value = "  -5.37%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.12"
$ws.Range("E46").Value = "  -1.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0510"
$ws.Range("E47").Value = "  -1.72%  "

# Row 48
$ws.Range("D48").Value = "1.993.49"
$ws.Range("E48").Value = "  -1.26%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0671"
$ws.Range("E49").Value = "  +8.08%  "

# Row 50
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  -0.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.57"
$ws.Range("E51").Value = "  -4.91%  "
